# Add a new row (9) to the "data" sheet for "Université d'Aix Marseille",
# matching the pattern of the previously-appended rows (7, 8): plain values,
# no custom row styling, and column G left empty (no hyperlink cell yet for
# this entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Université d'Aix Marseille"
$ws.Range("C9").Value = "U"
$ws.Range("D9").Value = "Marseille"
$ws.Range("E9").Value = "58, bd Charles Livon 13284 Marseille Cedex 07"

# Column F ("code" column) reuses the same "7764" text value already used by
# rows 2, 4 and 8. Assigning the literal string directly would make Excel
# auto-coerce it to a number, so instead copy the existing text cell — this
# keeps it as shared-string text (matching F2/F4/F8) without pulling in a new
# number-format style.
$ws.Range("F8").Copy($ws.Range("F9"))
